{"js": "// Search for the existing \"Test bestand\" text (it is split across two\n// runs around a stale \"_GoBack\" bookmark) and replace it with a single,\n// clean run containing the same text, then append the new content.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"Test bestand\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"Test bestand\", \"Replace\");\n}\n\n// Drop the leftover \"_GoBack\" bookmark that used to sit inside the text.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Append a blank paragraph followed by the new sentence.\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nlet lastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nlastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Dit is een tekst geschreven door klaas\", \"After\");\nawait context.sync();\n", "ps1": "# Search-and-replace helper: find a target string in the document body\n# and replace it with the given replacement text.\nfunction Replace-DocText {\n    param($Document, $FindText, $ReplaceText)\n\n    $find = $Document.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    return $find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n}\n\n$d = $word.ActiveDocument\n\n# \"Test besta\" + a stray \"_GoBack\" bookmark + \"nd\" all render as\n# \"Test bestand\" -- normalize it into one clean run of text.\nReplace-DocText $d \"Test bestand\" \"Test bestand\" | Out-Null\n\n# Append a blank paragraph and a new sentence at the end of the document.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Dit is een tekst geschreven door klaas\"\n"}
